$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Menu items (column A) and prices (column B), rows 2-14
$items = @(
    @("Rice and Chicken", 20),
    @("Chicken soup", 10),
    @("Crispy corn", 10),
    @("Shrimp", 16),
    @("Chicken salad", 10),
    @("Grilled Chicken and Potatoes", 15),
    @("Chocolate cake", 5),
    @("Fresh Juice", 2),
    @("Ice cream", 5),
    @("Coffee", 2),
    @("Tea", 2),
    @("Wine", 10),
    @("Apple pie", 5)
)

$row = 2
foreach ($item in $items) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row++
}

# Column A width to fit the longer item names
$ws.Columns.Item(1).ColumnWidth = 20.5

# Selection matching the final state
$ws.Range("E8").Select()
